$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells whose new value would otherwise be coerced to a number
# (losing a significant trailing zero) are pre-formatted as text so the
# written value keeps its exact original string shape, matching the source data.
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"

$ws.Range("D2").Value = "65.712.75"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "3.486.12"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "579.73"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "160.84"
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.483.68"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("E9").Value = "  +3.29%  "
$ws.Range("E10").Value = "  -3.94%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("D13").Value = "4.086.68"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("E15").Value = "  -2.09%  "
$ws.Range("D16").Value = "28.78"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").Value = "65.693.56"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").Value = "3.538.85"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").Value = "6.42"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("D21").Value = "391.30"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").Value = "8.25"
$ws.Range("E22").Value = "  -3.98%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "73.59"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").Value = "9.71"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "6.44"
$ws.Range("E30").Value = "  +6.41%  "
$ws.Range("E31").Value = "  +3.48%  "
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").Value = "6.49"
$ws.Range("E34").Value = "  -4.18%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +1.11%  "
$ws.Range("E37").Value = "  +4.39%  "
$ws.Range("D38").Value = "163.11"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("E39").Value = "  +4.37%  "
$ws.Range("D40").Value = "3.080.91"
$ws.Range("E40").Value = "  +5.75%  "
$ws.Range("D41").Value = "0.0771"
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("D42").Value = "27.20"
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("E43").Value = "  -1.94%  "
$ws.Range("E44").Value = "  +1.34%  "
$ws.Range("D45").Value = "42.78"
$ws.Range("D46").Value = "0.776"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("D47").Value = "25.86"
$ws.Range("E47").Value = "  +8.23%  "
$ws.Range("E48").Value = "  +2.11%  "
$ws.Range("D49").Value = "2.24"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").Value = "312.12"
$ws.Range("E50").Value = "  +4.85%  "
$ws.Range("E51").Value = "  +1.46%  "
